# repull data, push all data, mean calculation
# Update the dSF (column F) values for specific rows to reflect the repulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    8  = 2
    15 = 0
    16 = 2
    20 = -2
    24 = 1
    26 = -1
    27 = -8
    29 = 0
    31 = 0
    33 = 1
    34 = 0
    45 = -1
    48 = -6
    50 = -2
    57 = -3
    61 = 3
    66 = -2
    70 = 0
    71 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
